# Load screen Code Updated
# Adds a new "Loads_TC004" test-data row to the CustomizeGrid, Add Load and
# View Load sheets, tweaks a couple of saved cell-selections, and moves the
# active tab / active sheet selection over to "View Load".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CustomizeGrid (sheet1): append rows 6 & 7 for Loads_TC004
# ---------------------------------------------------------------------
$wsGrid = $wb.Worksheets.Item("CustomizeGrid")

$wsGrid.Range("A6").Value = "Loads_TC004"
$wsGrid.Range("B6").Value = "Load Date;Shipper;Shipper Contact;Carrier Name;Status;Origin;Destination;Rate;Rate UOM;Commodity;Ready to Submit"
$wsGrid.Range("C6").Value = "SELECT"
$wsGrid.Range("D6").Value = "Webtable customized successfully"

$wsGrid.Range("A7").Value = "Loads_TC004"
$wsGrid.Range("B7").Value = "Invoice #;Status"
$wsGrid.Range("C7").Value = "SELECT"
$wsGrid.Range("D7").Value = "Webtable customized successfully"

$null = $wsGrid.Range("A7").Select()

# ---------------------------------------------------------------------
# Add Load (sheet2): append row 5 for Loads_TC004
# ---------------------------------------------------------------------
$wsAdd = $wb.Worksheets.Item("Add Load")

$wsAdd.Range("A5").Value = "Loads_TC004"
$wsAdd.Range("B5").Value = "New Day"
$wsAdd.Range("C5").Value = "Current Date"
$wsAdd.Range("D5").Value = "CP Shipper"
$wsAdd.Range("E5").Value = "TestContact"
$wsAdd.Range("F5").Value = "Corn"
$wsAdd.Range("G5").NumberFormat = "@"
$wsAdd.Range("G5").Value = "0.25"
$wsAdd.Range("H5").Value = "Bushels"
$wsAdd.Range("I5").Value = "Alaska"
$wsAdd.Range("J5").Value = "Roger"
$wsAdd.Range("K5").Value = "Added new load successfully"

$null = $wsAdd.Range("B5:J5").Select()

# ---------------------------------------------------------------------
# Edit Load (sheet3): just move the saved selection
# ---------------------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("Edit Load")
$null = $wsEdit.Range("A12").Select()

# ---------------------------------------------------------------------
# View Load (sheet4): append row 8 for Loads_TC004 and make this the
# active tab / active sheet
# ---------------------------------------------------------------------
$wsView = $wb.Worksheets.Item("View Load")

$wsView.Range("A8").Value = "Loads_TC004"
$wsView.Range("B8").Value = "Current Date"
$wsView.Range("C8").Value = "CP Shipper"
$wsView.Range("D8").Value = "TestContact"
$wsView.Range("E8").Value = "New Day"
$wsView.Range("F8").Value = "Open"
$wsView.Range("G8").Value = "Alaska"
$wsView.Range("H8").Value = "Roger"
$wsView.Range("I8").NumberFormat = "@"
$wsView.Range("I8").Value = "0.25"
$wsView.Range("J8").Value = "Bushels"
$wsView.Range("K8").Value = "Corn"
$wsView.Range("L8").Value = "NA"
$wsView.Range("M8").Value = "ADD"
$wsView.Range("N8").Value = "Webtable validated successfully"

$null = $wsView.Activate()
$null = $wsView.Range("A8").Select()
